$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 25.81602566666667
$ws.Range("H2").Value = 77.44807700000001
$ws.Range("I2").Value = 0.7742517153725241
$ws.Range("J2").Value = 0.7742517153725241
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 88.37814633333333
$ws.Range("N2").Value = 265.134439
$ws.Range("O2").Value = 0.7138016014383547
$ws.Range("P2").Value = 0.7138016014383547
$ws.Range("Q2").Value = 2281.572494113756
$ws.Range("R2").Value = 20534.1524470238
$ws.Range("S2").Value = 0.5526621143493009
$ws.Range("T2").Value = 0.5526621143493009

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 25.81602566666667
$ws.Range("H3").Value = 77.44807700000001
$ws.Range("I3").Value = 0.7742517153725241
$ws.Range("J3").Value = 0.7742517153725241
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 8.866675333333333
$ws.Range("N3").Value = 26.600026
$ws.Range("O3").Value = 0.07161325864989523
$ws.Range("P3").Value = 0.07161325864989525
$ws.Range("Q3").Value = 228.9023179833336
$ws.Range("R3").Value = 2060.120861850002
$ws.Range("S3").Value = 0.05544668835309764
$ws.Range("T3").Value = 0.05544668835309765

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 25.81602566666667
$ws.Range("H4").Value = 77.44807700000001
$ws.Range("I4").Value = 0.7742517153725241
$ws.Range("J4").Value = 0.7742517153725241
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 26.56849866666667
$ws.Range("N4").Value = 79.705496
$ws.Range("O4").Value = 0.2145851399117501
$ws.Range("P4").Value = 0.2145851399117501
$ws.Range("Q4").Value = 685.893043503466
$ws.Range("R4").Value = 6173.037391531193
$ws.Range("S4").Value = 0.1661429126701256
$ws.Range("T4").Value = 0.1661429126701256

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.718527666666667
$ws.Range("H5").Value = 8.155583
$ws.Range("I5").Value = 0.08153170965901445
$ws.Range("J5").Value = 0.08153170965901445
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 88.37814633333333
$ws.Range("N5").Value = 265.134439
$ws.Range("O5").Value = 0.7138016014383547
$ws.Range("P5").Value = 0.7138016014383547
$ws.Range("Q5").Value = 240.2584359358819
$ws.Range("R5").Value = 2162.325923422937
$ws.Range("S5").Value = 0.05819746492261148
$ws.Range("T5").Value = 0.05819746492261148

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.718527666666667
$ws.Range("H6").Value = 8.155583
$ws.Range("I6").Value = 0.08153170965901445
$ws.Range("J6").Value = 0.08153170965901445
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.866675333333333
$ws.Range("N6").Value = 26.600026
$ws.Range("O6").Value = 0.07161325864989523
$ws.Range("P6").Value = 0.07161325864989525
$ws.Range("Q6").Value = 24.10430220501755
$ws.Range("R6").Value = 216.938719845158
$ws.Range("S6").Value = 0.005838751411979163
$ws.Range("T6").Value = 0.005838751411979164

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.718527666666667
$ws.Range("H7").Value = 8.155583
$ws.Range("I7").Value = 0.08153170965901445
$ws.Range("J7").Value = 0.08153170965901445
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 26.56849866666667
$ws.Range("N7").Value = 79.705496
$ws.Range("O7").Value = 0.2145851399117501
$ws.Range("P7").Value = 0.2145851399117501
$ws.Range("Q7").Value = 72.22719868712977
$ws.Range("R7").Value = 650.0447881841679
$ws.Range("S7").Value = 0.0174954933244238
$ws.Range("T7").Value = 0.0174954933244238

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 4.808641333333334
$ws.Range("H8").Value = 14.425924
$ws.Range("I8").Value = 0.1442165749684613
$ws.Range("J8").Value = 0.1442165749684613
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 88.37814633333333
$ws.Range("N8").Value = 265.134439
$ws.Range("O8").Value = 0.7138016014383547
$ws.Range("P8").Value = 0.7138016014383547
$ws.Range("Q8").Value = 424.9788074218485
$ws.Range("R8").Value = 3824.809266796636
$ws.Range("S8").Value = 0.1029420221664422
$ws.Range("T8").Value = 0.1029420221664422

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 4.808641333333334
$ws.Range("H9").Value = 14.425924
$ws.Range("I9").Value = 0.1442165749684613
$ws.Range("J9").Value = 0.1442165749684613
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 8.866675333333333
$ws.Range("N9").Value = 26.600026
$ws.Range("O9").Value = 0.07161325864989523
$ws.Range("P9").Value = 0.07161325864989525
$ws.Range("Q9").Value = 42.63666149711378
$ws.Range("R9").Value = 383.7299534740241
$ws.Range("S9").Value = 0.01032781888481843
$ws.Range("T9").Value = 0.01032781888481843

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 4.808641333333334
$ws.Range("H10").Value = 14.425924
$ws.Range("I10").Value = 0.1442165749684613
$ws.Range("J10").Value = 0.1442165749684613
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 26.56849866666667
$ws.Range("N10").Value = 79.705496
$ws.Range("O10").Value = 0.2145851399117501
$ws.Range("P10").Value = 0.2145851399117501
$ws.Range("Q10").Value = 127.7583808531449
$ws.Range("R10").Value = 1149.825427678304
$ws.Range("S10").Value = 0.03094673391720066
$ws.Range("T10").Value = 0.03094673391720066
